# TestCargoDetails.xlsx - CSV Parsing Utils TestCases -> bumped up LineCoverage to 90%
#
# Updates a handful of CountryCode / ChargeableUnit / NetWeightUnit /
# VolumetricWeightUnit / PacksType sample values on the CargoDetails sheet,
# widens column C, and refreshes the sheet's view state (zoom + selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CargoDetails")

# --- Cell value edits (row 2 / row 3 sample data) -------------------------
$ws.Range("D2").Value  = "BBG"   # PacksType
$ws.Range("V2").Value  = "HEG"   # ChargeableUnit
$ws.Range("X2").Value  = "IND"   # CountryCode
$ws.Range("AU2").Value = "CG"    # VolumetricWeightUnit

$ws.Range("T3").Value  = "DT"    # NetWeightUnit
$ws.Range("V3").Value  = "HEG"   # ChargeableUnit
$ws.Range("X3").Value  = "IND"   # CountryCode
$ws.Range("AU3").Value = "CG"    # VolumetricWeightUnit

# --- Column width (column C widened, no longer auto bestFit) -------------
$ws.Columns.Item(3).ColumnWidth = 11

# --- View state: zoom in and move the selection -------------------------
$win = $excel.ActiveWindow
$win.Zoom = 211
$win.ScrollRow = 1
$win.ScrollColumn = 11
$ws.Range("V4").Select()
